$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 140
$ws.Range("F3").Value = 960
$ws.Range("F4").Value = 599
$ws.Range("F5").Value = 2916
$ws.Range("F6").Value = 2916
$ws.Range("F7").Value = 785
$ws.Range("F8").Value = 591
$ws.Range("F10").Value = 406
$ws.Range("F11").Value = 655
$ws.Range("F12").Value = 377
$ws.Range("F13").Value = 490
$ws.Range("F14").Value = 525
$ws.Range("F15").Value = 2151
$ws.Range("F16").Value = 1258
$ws.Range("F17").Value = 740
$ws.Range("F19").Value = 15
$ws.Range("F20").Value = 2662
$ws.Range("F22").Value = 37
$ws.Range("F23").Value = 1048
$ws.Range("F24").Value = 528
$ws.Range("F25").Value = 513
$ws.Range("F26").Value = 592
$ws.Range("F27").Value = 593
$ws.Range("F28").Value = 7
$ws.Range("F29").Value = 20
$ws.Range("F31").Value = 552
$ws.Range("F32").Value = 564
$ws.Range("F34").Value = 113
$ws.Range("F35").Value = 386
$ws.Range("F36").Value = 4661
$ws.Range("F37").Value = 237
$ws.Range("F38").Value = 10

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 18
$ws.Range("F4").Value = 1
$ws.Range("F5").Value = 69
$ws.Range("F6").Value = 8
$ws.Range("F25").Value = 294
$ws.Range("F27").Value = 86
$ws.Range("F31").Value = 15
$ws.Range("F36").Value = 521
$ws.Range("F37").Value = 521
$ws.Range("F38").Value = 14

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 1462
$ws.Range("F5").Value = 566
$ws.Range("F6").Value = 225
$ws.Range("F7").Value = 240

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 1462
$ws.Range("F3").Value = 566
$ws.Range("F4").Value = 140
$ws.Range("F6").Value = 225
$ws.Range("F7").Value = 18
$ws.Range("F8").Value = 960
$ws.Range("F9").Value = 599
$ws.Range("F10").Value = 2916
$ws.Range("F11").Value = 2916
$ws.Range("F12").Value = 785
$ws.Range("F13").Value = 591
$ws.Range("F15").Value = 406
$ws.Range("F16").Value = 655
$ws.Range("F17").Value = 2
$ws.Range("F18").Value = 8
$ws.Range("F19").Value = 490
$ws.Range("F24").Value = 525
$ws.Range("F25").Value = 2151
$ws.Range("F26").Value = 740
$ws.Range("F28").Value = 2662
$ws.Range("F30").Value = 1048
$ws.Range("F31").Value = 528
$ws.Range("F33").Value = 240
$ws.Range("F37").Value = 513
$ws.Range("F38").Value = 593
$ws.Range("F39").Value = 593
$ws.Range("F41").Value = 552
$ws.Range("F42").Value = 564
$ws.Range("F43").Value = 294
$ws.Range("F45").Value = 113
$ws.Range("F46").Value = 386
$ws.Range("F48").Value = 4661
$ws.Range("F50").Value = 521
$ws.Range("F51").Value = 14
